$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7193536307022441
$ws.Range("C2").Value = 0.1900541378288949
$ws.Range("D2").Value = 0.01839289231952534
$ws.Range("F2").Value = 0.3904829946890942
$ws.Range("G2").Value = 0.002372792662579092
$ws.Range("I2").Value = 0.262931302456419
$ws.Range("N2").Value = 0.8448740361739766
$ws.Range("O2").Value = 1.209995590113039

$ws.Range("B3").Value = 0.6297058750278666
$ws.Range("C3").Value = 0.1676240295245748
$ws.Range("D3").Value = 0.01613704970889529
$ws.Range("F3").Value = 0.386585627643889
$ws.Range("G3").Value = 0.002375102843149508
$ws.Range("I3").Value = 0.2674990674277549
$ws.Range("N3").Value = 0.8395351397227842
$ws.Range("O3").Value = 1.209249759535481

$ws.Range("B4").Value = 0.5745020706209516
$ws.Range("C4").Value = 0.1537789512613017
$ws.Range("D4").Value = 0.01474554807747808
$ws.Range("F4").Value = 0.3845332939323782
$ws.Range("G4").Value = 0.002376597830151557
$ws.Range("I4").Value = 0.2705357709640523
$ws.Range("N4").Value = 0.8365931547189689
$ws.Range("O4").Value = 1.209874358393577

$ws.Range("B5").Value = 0.5519675264974637
$ws.Range("C5").Value = 0.1481190086781226
$ws.Range("D5").Value = 0.01417692723481423
$ws.Range("F5").Value = 0.383782481684193
$ws.Range("G5").Value = 0.00237722634860628
$ws.Range("I5").Value = 0.2718314826330115
$ws.Range("N5").Value = 0.8354791436990041
$ws.Range("O5").Value = 1.210400557845276

$ws.Range("B6").Value = 0.5482234021986017
$ws.Range("C6").Value = 0.1471781056609984
$ws.Range("D6").Value = 0.01408241431169444
$ws.Range("F6").Value = 0.3836629718092013
$ws.Range("G6").Value = 0.002377331880862537
$ws.Range("I6").Value = 0.2720501477939354
$ws.Range("N6").Value = 0.8352992995037027
$ws.Range("O6").Value = 1.210504325120752

$ws.Range("B7").Value = 0.5741983158034429
$ws.Range("C7").Value = 0.1537026915267177
$ws.Range("D7").Value = 0.01473788577484214
$ws.Range("F7").Value = 0.3845228220829924
$ws.Range("G7").Value = 0.002376606228410459
$ws.Range("I7").Value = 0.2705530097816204
$ws.Range("N7").Value = 0.8365777866644777
$ws.Range("O7").Value = 1.209880355643492

$ws.Range("B8").Value = 0.6884771769010172
$ws.Range("C8").Value = 0.1823355828501576
$ws.Range("D8").Value = 0.01761642931087692
$ws.Range("F8").Value = 0.3890683953646601
$ws.Range("G8").Value = 0.002373573363362596
$ws.Range("I8").Value = 0.2644580240362941
$ws.Range("N8").Value = 0.8429636292562321
$ws.Range("O8").Value = 1.209513415372641

$ws.Range("B9").Value = 0.9112498146081975
$ws.Range("C9").Value = 0.2378924128958317
$ws.Range("D9").Value = 0.02320900784459923
$ws.Range("F9").Value = 0.4006923993217981
$ws.Range("G9").Value = 0.002368230583590258
$ws.Range("I9").Value = 0.2543528256984935
$ws.Range("N9").Value = 0.8581394008761691
$ws.Range("O9").Value = 1.21740997196946

$ws.Range("B10").Value = 1.074045439065003
$ws.Range("C10").Value = 0.2783342793482575
$ws.Range("D10").Value = 0.02728450423169448
$ws.Range("F10").Value = 0.4108962254752626
$ws.Range("G10").Value = 0.002364670285959984
$ws.Range("I10").Value = 0.2480621413209931
$ws.Range("N10").Value = 0.8708903624068967
$ws.Range("O10").Value = 1.228504113608011

$ws.Range("B11").Value = 1.147902262880621
$ws.Range("C11").Value = 0.2966478488964697
$ws.Range("D11").Value = 0.02913100488108711
$ws.Range("F11").Value = 0.4159020043155266
$ws.Range("G11").Value = 0.002363129120458771
$ws.Range("I11").Value = 0.2454481151621124
$ws.Range("N11").Value = 0.877035636054444
$ws.Range("O11").Value = 1.234709061950497

$ws.Range("B12").Value = 1.175839758147958
$ws.Range("C12").Value = 0.3035703652176096
$ws.Range("D12").Value = 0.02982911932497956
$ws.Range("F12").Value = 0.4178500736324509
$ws.Range("G12").Value = 0.002362556742201472
$ws.Range("I12").Value = 0.244493997936214
$ws.Range("N12").Value = 0.8794119721224547
$ws.Range("O12").Value = 1.237225892906537

$ws.Range("B13").Value = 1.169824302329857
$ws.Range("C13").Value = 0.3020800370059078
$ws.Range("D13").Value = 0.02967881816118734
$ws.Range("F13").Value = 0.4174281852672266
$ws.Range("G13").Value = 0.002362679515345047
$ws.Range("I13").Value = 0.2446978914717732
$ws.Range("N13").Value = 0.8788980001892384
$ws.Range("O13").Value = 1.236676404578731

$ws.Range("B14").Value = 1.150201318390089
$ws.Range("C14").Value = 0.2972176200306933
$ws.Range("D14").Value = 0.02918846176920198
$ws.Range("F14").Value = 0.4160612203556582
$ws.Range("G14").Value = 0.002363081805743104
$ws.Range("I14").Value = 0.2453689018813812
$ws.Range("N14").Value = 0.8772301533792302
$ws.Range("O14").Value = 1.234912769712281

$ws.Range("B15").Value = 1.138177652391676
$ws.Range("C15").Value = 0.2942376177350638
$ws.Range("D15").Value = 0.02888795782956066
$ws.Range("F15").Value = 0.4152307548487144
$ws.Range("G15").Value = 0.002363329680937809
$ws.Range("I15").Value = 0.2457845762195134
$ws.Range("N15").Value = 0.8762149537913615
$ws.Range("O15").Value = 1.233854278486319

$ws.Range("B16").Value = 1.069214552372614
$ws.Range("C16").Value = 0.2771357266051382
$ws.Range("D16").Value = 0.02716367729719593
$ws.Range("F16").Value = 0.4105764233163356
$ws.Range("G16").Value = 0.002364772578738683
$ws.Range("I16").Value = 0.2482379698281569
$ws.Range("N16").Value = 0.8704956644467217
$ws.Range("O16").Value = 1.228121963081236

$ws.Range("B17").Value = 1.026855516245973
$ws.Range("C17").Value = 0.2666225580576622
$ws.Range("D17").Value = 0.02610394560927887
$ws.Range("F17").Value = 0.4078144770082375
$ws.Range("G17").Value = 0.002365677801805767
$ws.Range("I17").Value = 0.249806578651409
$ws.Range("N17").Value = 0.8670751421911547
$ws.Range("O17").Value = 1.224902442492521

$ws.Range("B18").Value = 1.002473025245479
$ws.Range("C18").Value = 0.2605678134714537
$ws.Range("D18").Value = 0.02549371602437844
$ws.Range("F18").Value = 0.4062601320420143
$ws.Range("G18").Value = 0.002366205847610272
$ws.Range("I18").Value = 0.2507321010039405
$ws.Range("N18").Value = 0.8651402214437667
$ws.Range("O18").Value = 1.223159635577872

$ws.Range("B19").Value = 0.9942143757958206
$ws.Range("C19").Value = 0.2585164455115034
$ws.Range("D19").Value = 0.02528698383820682
$ws.Range("F19").Value = 0.4057397361402835
$ws.Range("G19").Value = 0.002366385905014952
$ws.Range("I19").Value = 0.2510494642621417
$ws.Range("N19").Value = 0.864490677876816
$ws.Range("O19").Value = 1.222588249236253

$ws.Range("B20").Value = 1.031366655457475
$ws.Range("C20").Value = 0.2677425180015973
$ws.Range("D20").Value = 0.02621682860331731
$ws.Range("F20").Value = 0.4081049447367775
$ws.Range("G20").Value = 0.002365580675365575
$ws.Range("I20").Value = 0.2496371850244827
$ws.Range("N20").Value = 0.8674359041569204
$ws.Range("O20").Value = 1.225233882904348

$ws.Range("B21").Value = 1.15596590570442
$ws.Range("C21").Value = 0.2986461704933561
$ws.Range("D21").Value = 0.02933252193187741
$ws.Range("F21").Value = 0.416461305385063
$ws.Range("G21").Value = 0.002362963339106333
$ws.Range("I21").Value = 0.2451708382295976
$ws.Range("N21").Value = 0.877718706250846
$ws.Range("O21").Value = 1.235426250751459

$ws.Range("B22").Value = 1.237220237256736
$ws.Range("C22").Value = 0.3187708243492864
$ws.Range("D22").Value = 0.03136228741462332
$ws.Range("F22").Value = 0.4222286852526764
$ws.Range("G22").Value = 0.002361318178786648
$ws.Range("I22").Value = 0.2424603139876353
$ws.Range("N22").Value = 0.8847260553541076
$ws.Range("O22").Value = 1.243062118870597

$ws.Range("B23").Value = 1.193870183305307
$ws.Range("C23").Value = 0.3080367016900141
$ws.Range("D23").Value = 0.03027957388298574
$ws.Range("F23").Value = 0.4191224780841125
$ws.Range("G23").Value = 0.002362190262442783
$ws.Range("I23").Value = 0.2438878461632257
$ws.Range("N23").Value = 0.8809599527064904
$ws.Range("O23").Value = 1.238897339760541

$ws.Range("B24").Value = 1.029327261850199
$ws.Range("C24").Value = 0.2672362170505096
$ws.Range("D24").Value = 0.02616579723769519
$ws.Range("F24").Value = 0.4079735198302998
$ws.Range("G24").Value = 0.002365624562674138
$ws.Range("I24").Value = 0.2497136940450275
$ws.Range("N24").Value = 0.8672727052681637
$ws.Range("O24").Value = 1.225083701898114

$ws.Range("B25").Value = 0.8511329009566566
$ws.Range("C25").Value = 0.2229277408067105
$ws.Range("D25").Value = 0.02170181727331055
$ws.Range("F25").Value = 0.3972564633422238
$ws.Range("G25").Value = 0.002369611584981067
$ws.Range("I25").Value = 0.2568880231475497
$ws.Range("N25").Value = 0.8537515999508685
$ws.Range("O25").Value = 1.214347185201547
